$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row (row 6): date in column A (same format as A5), text in column B
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Value = 42692

$ws.Range("B6").Value = "addUser nonEmtpy Validation added and some(user) handling from getUser function implemented"

# Update selection to match post-edit state
$ws.Range("B7").Select()
